# Generate Report for Archive
#
# 1. Update the "Status" value shown for the single localized file from
#    "Ready for handoff" to "In Translation" on all three sheets
#    (Overview!E2/F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "Status" column(s) on all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text change: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes ---
# Target stored width ~13.41 chars; this runtime quantizes ColumnWidth to
# 1/6-character pixel steps, so 12.5 is the value that lands on the closest
# achievable stored width (80/6 = 13.3333...).
# (Column indices used numerically - E=5, F=6, C=3 - to address them.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5

Write-Host "Applied status text + column width updates"
